$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ======================================================================
# 1) Rewrite Sheet1 rows 7-13 : old upstream/dwnstream pool-width block
#    is replaced by Municipality + tide-prediction fields.
# ======================================================================
$ws1.Range("A7").Value2 = "Municipality"
$ws1.Range("B7").Value2 = "Data Sheet - SITE"
$ws1.Range("C7").Value2 = "G11"
$ws1.Range("D7").Value2 = "Field"

$ws1.Range("A8").Value2 = "streamName"
$ws1.Range("B8").Value2 = "Data Sheet - SITE"
$ws1.Range("C8").Value2 = "G12"
$ws1.Range("D8").Value2 = "Field"

$ws1.Range("A9").Value2 = "roadName"
$ws1.Range("B9").Value2 = "Data Sheet - SITE"
$ws1.Range("C9").Value2 = "G13"
$ws1.Range("D9").Value2 = "Field"

$ws1.Range("A10").Value2 = "TidePredictTimeHigh"
$ws1.Range("B10").Value2 = "Data Sheet - SITE"
$ws1.Range("C10").Value2 = "AA13"
$ws1.Range("D10").Value2 = "Field"

$ws1.Range("A11").Value2 = "TidePredictElevationHigh"
$ws1.Range("B11").Value2 = "Data Sheet - SITE"
$ws1.Range("C11").Value2 = "AA14"
$ws1.Range("D11").Value2 = "Field"

$ws1.Range("A12").Value2 = "TidePredictTimeLow"
$ws1.Range("B12").Value2 = "Data Sheet - SITE"
$ws1.Range("C12").Value2 = "AE13"
$ws1.Range("D12").Value2 = "Field"

$ws1.Range("A13").Value2 = "TidePredictElevationLow"
$ws1.Range("B13").Value2 = "Data Sheet - SITE"
$ws1.Range("C13").Value2 = "AE14"
$ws1.Range("D13").Value2 = "Field"

# Rows 14-25 (CrossingType ... CrosDim_dwnD) are unchanged by the edit.

# ======================================================================
# 2) Insert the new tide-perch rows (26-29), then re-append the
#    channel-width / pool-width / LiDAR block (30-34, moved down from
#    its old 7-11 position), then the new headwall/scour/assessment
#    rows (35-42).
# ======================================================================
$ws1.Range("A26").Value2 = "LowTidePerch_upStream"
$ws1.Range("B26").Value2 = "Data Sheet - SITE"
$ws1.Range("C26").Value2 = "Z55"
$ws1.Range("D26").Value2 = "Field"

$ws1.Range("A27").Value2 = "LowTidePerch_dwnStream"
$ws1.Range("B27").Value2 = "Data Sheet - SITE"
$ws1.Range("C27").Value2 = "AD55"
$ws1.Range("D27").Value2 = "Field"

$ws1.Range("A28").Value2 = "HighTidePerch_upStream"
$ws1.Range("B28").Value2 = "Data Sheet - SITE"
$ws1.Range("C28").Value2 = "Z56"
$ws1.Range("D28").Value2 = "Field"

$ws1.Range("A29").Value2 = "HighTidePerch_dwnStream"
$ws1.Range("B29").Value2 = "Data Sheet - SITE"
$ws1.Range("C29").Value2 = "AD56"
$ws1.Range("D29").Value2 = "Field"

$ws1.Range("A30").Value2 = "upstreamChannelwidth"
$ws1.Range("B30").Value2 = "Data Sheet - SUMMARY"
$ws1.Range("C30").Value2 = "F243"
$ws1.Range("D30").Value2 = "Desktop"

$ws1.Range("A31").Value2 = "dwnstreamChannelwidth"
$ws1.Range("B31").Value2 = "Data Sheet - SITE"
$ws1.Range("C31").Value2 = "K243"
$ws1.Range("D31").Value2 = "Desktop"

$ws1.Range("A32").Value2 = "upstreammaxPoolwidth"
$ws1.Range("B32").Value2 = "Data Sheet - SITE"
$ws1.Range("C32").Value2 = "F245"
$ws1.Range("D32").Value2 = "Desktop"

$ws1.Range("A33").Value2 = "dwnstreammaxPoolwidth"
$ws1.Range("B33").Value2 = "Data Sheet - SUMMARY"
$ws1.Range("C33").Value2 = "K245"
$ws1.Range("D33").Value2 = "Desktop"

$ws1.Range("A34").Value2 = "LiDarHt_CL"
$ws1.Range("B34").Value2 = "Data Sheet - SUMMARY"
$ws1.Range("C34").Value2 = "J54"
$ws1.Range("D34").Value2 = "Desktop"

$ws1.Range("A35").Value2 = "HeadwallMaterial_upStream"
$ws1.Range("B35").Value2 = "Data Sheet - SITE"
$ws1.Range("C35").Value2 = "B61"
$ws1.Range("D35").Value2 = "Field"

$ws1.Range("A36").Value2 = "HeadwallCondition_upStream"
$ws1.Range("B36").Value2 = "Data Sheet - SITE"
$ws1.Range("C36").Value2 = "I62"
$ws1.Range("D36").Value2 = "Field"

$ws1.Range("A37").Value2 = "WindwallCondition_upStream"
$ws1.Range("B37").Value2 = "Data Sheet - SITE"
$ws1.Range("C37").Value2 = "U62"
$ws1.Range("D37").Value2 = "Field"

$ws1.Range("A38").Value2 = "ScourStructure_upStream"
$ws1.Range("B38").Value2 = "Data Sheet - SITE"
$ws1.Range("C38").Value2 = "Y61"
$ws1.Range("D38").Value2 = "Field"

$ws1.Range("A39").Value2 = "ScourSeverity_upStream"
$ws1.Range("B39").Value2 = "Data Sheet - SITE"
$ws1.Range("C39").Value2 = "AD62"
$ws1.Range("D39").Value2 = "Field"

$ws1.Range("A40").Value2 = "CrossingConditionEval"
$ws1.Range("B40").Value2 = "Data Sheet - SUMMARY"
$ws1.Range("C40").Value2 = "N12"
$ws1.Range("D40").Value2 = "Assessment"

$ws1.Range("A41").Value2 = "TidalRngRatio"
$ws1.Range("B41").Value2 = "Data Sheet - SUMMARY"
$ws1.Range("C41").Value2 = "N14"
$ws1.Range("D41").Value2 = "Assessment"

$ws1.Range("A42").Value2 = "GeneralAssessmentNotes"
$ws1.Range("B42").Value2 = "Data Sheet - SITE"

# ======================================================================
# 3) Add Sheet2 (the new "key" lookup sheet used by the Sheet/Cell
#    drop-down), positioned after Sheet1.
# ======================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet, 1, $null)
$ws2.Name = "Sheet2"

$ws2.Columns.Item(1).ColumnWidth = 20.75

$ws2.Range("A1").Value2 = "key"
$ws2.Range("A2").Value2 = "Data Sheet - SITE"
$ws2.Range("A3").Value2 = "Data Sheet - SUMMARY"

$ws2.Range("A2:A3").Validation.Add(3, 1, 1, "`$B`$4:`$B`$18")
$ws2.Range("A6").Select()

# ======================================================================
# 4) Point Sheet1's "Sheet" column validation at the new Sheet2 list
#    instead of the old in-sheet $B$4:$B$7 range.
# ======================================================================
$ws1.Range("B2:B1048576").Validation.Delete()
$ws1.Range("B2:B39").Validation.Add(3, 1, 1, "Sheet2!`$A`$2:`$A`$3")

# ======================================================================
# 5) Restore Sheet1 as the active sheet/tab and match the author's
#    final on-screen selection (whole row 40 selected).
# ======================================================================
$ws1.Activate()
$ws1.Rows(40).Select()
